# Allocation_Template.xlsx edit:
#  - add a new "Fund_Accounts" worksheet (FundingRequest data) as the last tab
#  - populate header row (bold) + 4 data rows
#  - mirror the workbook's existing page-setup conventions on the new sheet
#  - nudge a couple of cursor/selection positions on existing sheets
#    (Accounts, Tax_Status) to match where the author had last clicked
#  - leave Other_inputs' tab no longer "selected" (Fund_Accounts becomes
#    the active tab once we select a cell on it, last)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet after the last existing one (Tax_Status)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Fund_Accounts"

# ---------------------------------------------------------------------
# 2. Header row - bold, matching the workbook's existing bold-header style
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Funding Number"
$ws.Range("B1").Value = "Account to Fund"
$ws.Range("C1").Value = "Amount to Fund"
$ws.Range("A1:C1").Font.Bold = $true

# ---------------------------------------------------------------------
# 3. Data rows (FundingRequest records)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "A11"
$ws.Range("C2").Value = 1000

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "A3"
$ws.Range("C3").Value = 2000

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "A4"
$ws.Range("C4").Value = 2000

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "A2"
$ws.Range("C5").Value = 9999999999999

# ---------------------------------------------------------------------
# 4. Row heights - match the workbook's 12.8pt default row height
# ---------------------------------------------------------------------
$ws.Range("A1:C5").RowHeight = 12.8

# ---------------------------------------------------------------------
# 5. Column widths - approximate the template's layout
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.36
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 20.11
$ws.Columns.Item(4).ColumnWidth = 21.87

# ---------------------------------------------------------------------
# 6. Page setup - reuse the same margins / header / footer as the rest
#    of the workbook (Times New Roman title / page-number footer)
# ---------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7
$ps.CenterHeader = "&""Times New Roman,Regular""&12&A"
$ps.CenterFooter = "&""Times New Roman,Regular""&12Page &P"
$ps.PrintGridlines = $false
$ps.PrintHeadings = $false
$ps.Zoom = 100
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.Orientation = 1
$ps.PaperSize = 1
$ps.CenterHorizontally = $false
$ps.CenterVertically = $false

# ---------------------------------------------------------------------
# 7. Selections - match where the author's cursor ended up on each tab.
#    Selecting on Fund_Accounts last makes it the active tab/sheet.
# ---------------------------------------------------------------------
$wsAccounts = $wb.Worksheets.Item("Accounts")
[void]$wsAccounts.Range("Q22").Select()

$wsTaxStatus = $wb.Worksheets.Item("Tax_Status")
[void]$wsTaxStatus.Range("J12").Select()

[void]$ws.Range("B8").Select()

Write-Host "Fund_Accounts sheet added with FundingRequest data."
